$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (price, volume%, and a few coin name/link reorderings)
$updates = @(
    @{ Cell = "D2"; Value = "309.57" }
    @{ Cell = "E2"; Value = "-4.04%" }
    @{ Cell = "D3"; Value = "49.49" }
    @{ Cell = "E3"; Value = "2.14%" }
    @{ Cell = "D4"; Value = "5.133" }
    @{ Cell = "E4"; Value = "-2.76%" }
    @{ Cell = "D5"; Value = "0.07758" }
    @{ Cell = "E5"; Value = "-4.31%" }
    @{ Cell = "D6"; Value = "4.530" }
    @{ Cell = "E6"; Value = "-1.14%" }
    @{ Cell = "D7"; Value = "1.376" }
    @{ Cell = "E7"; Value = "14.60%" }
    @{ Cell = "D8"; Value = "1.572" }
    @{ Cell = "E8"; Value = "-4.25%" }
    @{ Cell = "D9"; Value = "0.1217" }
    @{ Cell = "E9"; Value = "-6.45%" }
    @{ Cell = "D10"; Value = "0.1982" }
    @{ Cell = "E10"; Value = "1.89%" }
    @{ Cell = "D11"; Value = "0.04733" }
    @{ Cell = "E11"; Value = "2.09%" }
    @{ Cell = "D12"; Value = "0.09373" }
    @{ Cell = "E12"; Value = "-1.71%" }
    @{ Cell = "E13"; Value = "-0.53%" }
    @{ Cell = "D14"; Value = "0.001254" }
    @{ Cell = "E14"; Value = "-5.52%" }
    @{ Cell = "B15"; Value = "TigerCash" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D15"; Value = "0.005795" }
    @{ Cell = "E15"; Value = "-0.65%" }
    @{ Cell = "B16"; Value = "UpBots" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt" }
    @{ Cell = "D16"; Value = "0.007509" }
    @{ Cell = "E16"; Value = "2,021.82%" }
    @{ Cell = "B17"; Value = "LEO" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D17"; Value = "3.336" }
    @{ Cell = "E17"; Value = "-0.21%" }
    @{ Cell = "B18"; Value = "BTSEToken" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Cell = "D18"; Value = "2.434" }
    @{ Cell = "E18"; Value = "0.20%" }
    @{ Cell = "B19"; Value = "BitpandaEcosystemToken" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" }
    @{ Cell = "D19"; Value = "0.3395" }
    @{ Cell = "E19"; Value = "-0.20%" }
    @{ Cell = "B20"; Value = "MCDex" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" }
    @{ Cell = "D20"; Value = "7.957" }
    @{ Cell = "E20"; Value = "-2.03%" }
    @{ Cell = "B21"; Value = "ProBitToken" }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob" }
    @{ Cell = "D21"; Value = "0.1358" }
    @{ Cell = "E21"; Value = "-2.97%" }
    @{ Cell = "B22"; Value = "ZBToken" }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb" }
    @{ Cell = "D22"; Value = "0.3033" }
    @{ Cell = "E22"; Value = "-2.94%" }
    @{ Cell = "B23"; Value = "CoinExToken" }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" }
    @{ Cell = "D23"; Value = "0.04168" }
    @{ Cell = "E23"; Value = "2.15%" }
    @{ Cell = "D24"; Value = "0.001270" }
    @{ Cell = "E24"; Value = "-2.83%" }
    @{ Cell = "D25"; Value = "0.003924" }
    @{ Cell = "E25"; Value = "-7.72%" }
    @{ Cell = "D26"; Value = "0.0001348" }
    @{ Cell = "E26"; Value = "-0.16%" }
    @{ Cell = "D38"; Value = "0.02603" }
    @{ Cell = "E38"; Value = "-3.92%" }
    @{ Cell = "D39"; Value = "0.06227" }
    @{ Cell = "E39"; Value = "9.58%" }
    @{ Cell = "D40"; Value = "0.01098" }
    @{ Cell = "E40"; Value = "74.34%" }
    @{ Cell = "D41"; Value = "0.007943" }
    @{ Cell = "E41"; Value = "3.34%" }
    @{ Cell = "D42"; Value = "0.1421" }
    @{ Cell = "E42"; Value = "-1.32%" }
    @{ Cell = "D43"; Value = "0.008369" }
    @{ Cell = "E43"; Value = "8.67%" }
    @{ Cell = "D44"; Value = "0.008323" }
    @{ Cell = "E44"; Value = "2.74%" }
    @{ Cell = "D45"; Value = "0.3125" }
    @{ Cell = "E45"; Value = "-2.15%" }
    @{ Cell = "D46"; Value = "0.00007626" }
    @{ Cell = "E46"; Value = "9.17%" }
    @{ Cell = "D47"; Value = "0.00000000749" }
    @{ Cell = "E47"; Value = "-0.14%" }
    @{ Cell = "B48"; Value = "BOLO" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo" }
    @{ Cell = "D48"; Value = "0.05319" }
    @{ Cell = "E48"; Value = "-14.76%" }
    @{ Cell = "B49"; Value = "CoinbaseStockToken" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin" }
    @{ Cell = "D49"; Value = "0.002616" }
    @{ Cell = "E49"; Value = "-34.59%" }
    @{ Cell = "D50"; Value = "0.00002097" }
    @{ Cell = "E50"; Value = "-0.14%" }
    @{ Cell = "D51"; Value = "0.0001997" }
    @{ Cell = "E51"; Value = "-0.14%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
